$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("C2").Value = 5214
$ws.Range("C13").Select()
